# Scheduled runner update: refresh market-board derived profit figures
# (currentAveragePrice / NQ / HQ / LevePrice / LeveProfit columns, H:N)
# for a batch of Leve rows across the ALC/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 332.3889
$ws.Range("I53").Value = 121.46154
$ws.Range("J53").Value = 880.8
$ws.Range("K53").Value = 121.46154
$ws.Range("L53").Value = 880.8
$ws.Range("M53").Value = 515.53846
$ws.Range("N53").Value = -2154.8

$ws.Range("H113").Value = 3523.4443
$ws.Range("J113").Value = 2386.5715
$ws.Range("L113").Value = 2386.5715
$ws.Range("N113").Value = -8894.5715

$ws.Range("H132").Value = 1928.6735
$ws.Range("I132").Value = 1589.6052
$ws.Range("J132").Value = 3100
$ws.Range("K132").Value = 4768.8156
$ws.Range("L132").Value = 9300
$ws.Range("M132").Value = -2238.8156
$ws.Range("N132").Value = -14360

$ws.Range("H137").Value = 7938365.5
$ws.Range("I137").Value = 1498.2439
$ws.Range("J137").Value = 22729800
$ws.Range("K137").Value = 4494.7317
$ws.Range("L137").Value = 68189400
$ws.Range("M137").Value = -1944.7317
$ws.Range("N137").Value = -68194500

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 450
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").Value = $null

$ws.Range("H64").Value = 907.1539
$ws.Range("I64").Value = 663.3333
$ws.Range("J64").Value = 980.3
$ws.Range("K64").Value = 663.3333
$ws.Range("L64").Value = 980.3
$ws.Range("M64").Value = -438.3333
$ws.Range("N64").Value = -1430.3

$ws.Range("H67").Value = 907.1539
$ws.Range("I67").Value = 663.3333
$ws.Range("J67").Value = 980.3
$ws.Range("K67").Value = 663.3333
$ws.Range("L67").Value = 980.3
$ws.Range("M67").Value = 116.6667
$ws.Range("N67").Value = -2540.3

$ws.Range("H81").Value = 19300
$ws.Range("J81").Value = 19300
$ws.Range("L81").Value = 19300
$ws.Range("N81").Value = -21422

$ws.Range("H84").Value = 19300
$ws.Range("J84").Value = 19300
$ws.Range("L84").Value = 57900
$ws.Range("N84").Value = -68508

$ws.Range("H134").Value = 71552.516
$ws.Range("I134").Value = 2887.4092
$ws.Range("J134").Value = 287357.16
$ws.Range("K134").Value = 8662.2276
$ws.Range("L134").Value = 862071.48
$ws.Range("M134").Value = -6127.2276
$ws.Range("N134").Value = -867141.48

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2254.0571
$ws.Range("I99").Value = 2124.2068
$ws.Range("J99").Value = 2881.6667
$ws.Range("K99").Value = 2124.2068
$ws.Range("L99").Value = 2881.6667
$ws.Range("M99").Value = -626.2067999999999
$ws.Range("N99").Value = -5877.6667

$ws.Range("H126").Value = 2254.0571
$ws.Range("I126").Value = 2124.2068
$ws.Range("J126").Value = 2881.6667
$ws.Range("K126").Value = 6372.6204
$ws.Range("L126").Value = 8645.000100000001
$ws.Range("M126").Value = -3902.6204
$ws.Range("N126").Value = -13585.0001

$ws.Range("H132").Value = 4983
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 5779.6
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 17338.8
$ws.Range("M132").Value = -470
$ws.Range("N132").Value = -22398.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 156775.64
$ws.Range("I107").Value = 334.15
$ws.Range("J107").Value = 261069.97
$ws.Range("K107").Value = 1002.45
$ws.Range("L107").Value = 783209.91
$ws.Range("M107").Value = 917.5500000000001
$ws.Range("N107").Value = -787049.91

$ws.Range("H131").Value = 14316367
$ws.Range("I131").Value = 41750460
$ws.Range("J131").Value = 2926.652
$ws.Range("K131").Value = 125251380
$ws.Range("L131").Value = 8779.956
$ws.Range("M131").Value = -125246340
$ws.Range("N131").Value = -18859.956

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4582.4
$ws.Range("I126").Value = 4637.3335
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 13912.0005
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = -11442.0005
$ws.Range("N126").Value = -18440

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1266.5
$ws.Range("I7").Value = 1299.6666
$ws.Range("J7").Value = 1255.4445
$ws.Range("K7").Value = 1299.6666
$ws.Range("L7").Value = 1255.4445
$ws.Range("M7").Value = -1187.6666
$ws.Range("N7").Value = -1479.4445

$ws.Range("H40").Value = 1158768.6
$ws.Range("I40").Value = 1300489.8
$ws.Range("K40").Value = 1300489.8
$ws.Range("M40").Value = -1300353.8

$ws.Range("H61").Value = 916.7273
$ws.Range("I61").Value = 916.7273
$ws.Range("K61").Value = 916.7273
$ws.Range("M61").Value = -714.7273

$ws.Range("H113").Value = 916.7273
$ws.Range("I113").Value = 916.7273
$ws.Range("K113").Value = 916.7273
$ws.Range("M113").Value = 1253.2727

$ws.Range("H126").Value = 1266.5
$ws.Range("I126").Value = 1299.6666
$ws.Range("J126").Value = 1255.4445
$ws.Range("K126").Value = 3898.9998
$ws.Range("L126").Value = 3766.3335
$ws.Range("M126").Value = -1428.9998
$ws.Range("N126").Value = -8706.333500000001

$ws.Range("H132").Value = 3033344.5
$ws.Range("I132").Value = 3499282.2
$ws.Range("K132").Value = 10497846.6
$ws.Range("M132").Value = -10495316.6

$ws.Range("H136").Value = 2235.5686
$ws.Range("I136").Value = 1361.0605
$ws.Range("J136").Value = 3838.8333
$ws.Range("K136").Value = 4083.1815
$ws.Range("L136").Value = 11516.4999
$ws.Range("M136").Value = -1533.1815
$ws.Range("N136").Value = -16616.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1163.3334
$ws.Range("I113").Value = 1250
$ws.Range("J113").Value = 990
$ws.Range("K113").Value = 3750
$ws.Range("L113").Value = 2970
$ws.Range("M113").Value = -1580
$ws.Range("N113").Value = -7310

$ws.Range("H122").Value = 2126.75
$ws.Range("I122").Value = 2102
$ws.Range("J122").Value = 2135
$ws.Range("K122").Value = 6306
$ws.Range("L122").Value = 6405
$ws.Range("M122").Value = -3856
$ws.Range("N122").Value = -11305

$ws.Range("H126").Value = 992.6667
$ws.Range("I126").Value = 468.8
$ws.Range("J126").Value = 1647.5
$ws.Range("K126").Value = 1406.4
$ws.Range("L126").Value = 4942.5
$ws.Range("M126").Value = 1063.6
$ws.Range("N126").Value = -9882.5

$ws.Range("H132").Value = 2734.0667
$ws.Range("I132").Value = 1038
$ws.Range("J132").Value = 3582.1
$ws.Range("K132").Value = 3114
$ws.Range("L132").Value = 10746.3
$ws.Range("M132").Value = -584
$ws.Range("N132").Value = -15806.3

$ws.Range("H136").Value = 2474.2744
$ws.Range("I136").Value = 2719.6128
$ws.Range("J136").Value = 2094
$ws.Range("K136").Value = 8158.8384
$ws.Range("L136").Value = 6282
$ws.Range("M136").Value = -5608.8384
$ws.Range("N136").Value = -11382
